$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Argentina): update 2026f (G) and 2027f (H)
$ws.Range("G2").Value = 14.5
$ws.Range("H2").Value = 14.3

# Row 5 (Colombia): update 2026f (G)
$ws.Range("G5").Value = 5.1

# Row 6 (Costa Rica): remove 2023 (D) value entirely
$ws.Range("D6").ClearContents()

# Row 12 (Mexico): add the full data series (2022-2027f)
$ws.Range("C12").Value = 29
$ws.Range("D12").Value = 27.9
$ws.Range("E12").Value = 25.5
$ws.Range("F12").Value = 25.2
$ws.Range("G12").Value = 24.9
$ws.Range("H12").Value = 24.6
